$d = $word.ActiveDocument

# 1. Remove the "_GoBack" bookmark that trails the "Important Links" heading.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# The document ends with three empty paragraphs. The first two become the
# new "Excel" section (heading + link); the third stays untouched as the
# document's trailing blank paragraph, and a brand-new paragraph is
# inserted before it to hold the relocated "_GoBack" bookmark.
$trailingCount = 3
$total = $d.Paragraphs.Count
$firstEmptyIndex = $total - $trailingCount + 1

# 2. First trailing empty paragraph -> bold "Excel" heading with a tab stop
#    and a trailing tab character.
$excelPara = $d.Paragraphs.Item($firstEmptyIndex)
$excelPara.Range.Text = "Excel"
$excelPara.Range.Font.Bold = $true
$excelPara.TabStops.Add(124.6)
$tabRange = $d.Range($excelPara.Range.End - 1, $excelPara.Range.End - 1)
$tabRange.InsertAfter([char]9)
$tabRange.Font.Bold = $true

# 3. Second trailing empty paragraph -> the Excel functions URL.
$urlPara = $d.Paragraphs.Item($firstEmptyIndex + 1)
$urlPara.Range.Text = "http://www.excelfunctions.net/Excel-Text-Functions.html"

# 4. Insert a new paragraph right after the URL paragraph (i.e. before the
#    final, still-empty trailing paragraph) and re-create the "_GoBack"
#    bookmark inside it so the paragraph ends up empty but bookmarked.
$urlPara.Range.InsertParagraphAfter()
$bmIndex = $firstEmptyIndex + 2
$bmPara = $d.Paragraphs.Item($bmIndex)
$bmPara.Range.Text = "X"
$bmRange = $d.Range($bmPara.Range.Start, $bmPara.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$bmPara = $d.Paragraphs.Item($bmIndex)
$clearRange = $d.Range($bmPara.Range.Start, $bmPara.Range.End - 1)
$clearRange.Text = ""
